# Insert a new week of "Acelga" price data (date 2023-01-13, serial 44939)
# ahead of the existing rows, pushing the remaining rows down by two and
# letting the oldest week's "Extra" row fall off and be recreated at the
# bottom of the block (rows 1066-1068 end up as Extra/Primera/Segunda for
# serial 44335, same as the original rows 1064-1066).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 1018; everything at/after 1018 shifts down by 2.
$ws.Rows.Item(1018).Insert()
$ws.Rows.Item(1018).Insert()

# New row 1018: Acelga "Primera" quality entry for the new week.
$ws.Range("A1018").Value = 6
$ws.Range("B1018").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1018").Value = "Metropolitana"
$ws.Range("D1018").Value = 44939
$ws.Range("E1018").Value = 13
$ws.Range("F1018").Value = 100112009
$ws.Range("G1018").Value = "Acelga"
$ws.Range("H1018").Value = "Sin especificar"
$ws.Range("I1018").Value = "Primera"
$ws.Range("J1018").Value = 440
$ws.Range("K1018").Value = 12000
$ws.Range("L1018").Value = 13000
$ws.Range("M1018").Value = 12432
$ws.Range("N1018").Value = "$/docena de atados"
$ws.Range("O1018").Value = "Región Metropolitana"
$ws.Range("P1018").Value = 4144
$ws.Range("Q1018").Value = 3
$ws.Range("R1018").Value = "Hortaliza"

# New row 1019: Acelga "Segunda" quality entry for the new week.
$ws.Range("A1019").Value = 6
$ws.Range("B1019").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1019").Value = "Metropolitana"
$ws.Range("D1019").Value = 44939
$ws.Range("E1019").Value = 13
$ws.Range("F1019").Value = 100112009
$ws.Range("G1019").Value = "Acelga"
$ws.Range("H1019").Value = "Sin especificar"
$ws.Range("I1019").Value = "Segunda"
$ws.Range("J1019").Value = 150
$ws.Range("K1019").Value = 10000
$ws.Range("L1019").Value = 10000
$ws.Range("M1019").Value = 10000
$ws.Range("N1019").Value = "$/docena de atados"
$ws.Range("O1019").Value = "Región Metropolitana"
$ws.Range("P1019").Value = 3333
$ws.Range("Q1019").Value = 3
$ws.Range("R1019").Value = "Hortaliza"
